$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (Late), pushing Late/Paid Date/Amount
# one column to the right (N -> O -> P -> Q), to make room for a new
# "Variable Instalments" related column as part of the Loan RBI change.
$ws.Columns("N").Insert()

# Inherit the column width from the column to the left (M), matching Excel's
# default behaviour when a column is inserted.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab, with P6 selected.
$ws.Activate() | Out-Null
$ws.Range("P6").Select() | Out-Null
